# Apply the "2/27 format changes on User CRM Pipeline feature file and
# step definitions" edit: append a new "Inventory" module row to the
# "modules" worksheet and move the active selection onto it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("modules")

# Make sure we're working on/viewing the "modules" sheet.
$ws.Activate()
$ws.Select()

# Add the new module entry in the first free row (A5).
$ws.Range("A5").Value = "Inventory"

# Move the active cell/selection to the newly written cell.
$ws.Range("A5").Select()
